$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.334.43'
$ws.Range("E2").Value = '  -1.86%  '
$ws.Range("D3").Value = '2.639.00'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("E4").Value = '  +0.05%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '581.29'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E5").Value = '  -2.16%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '156.24'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E6").Value = '  +0.37%  '
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.649'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E7").Value = '  +3.83%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  -3.10%  '
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("E12").Value = '  +0.07%  '
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '28.77'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("E14").Value = '  -4.15%  '
$ws.Range("D15").Value = '3.115.39'
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").Value = '64.172.90'
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("D17").Value = '2.636.81'
$ws.Range("E17").Value = '  +0.26%  '
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '12.23'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E18").Value = '  -2.54%  '
$ws.Range("E19").Value = '  -0.72%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '7.59'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E20").Value = '  +2.22%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '347.95'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  -0.14%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '67.88'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E23").Value = '  -1.20%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '1.76'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E24").Value = '  +7.71%  '
$ws.Range("E25").Value = '  -1.68%  '
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '9.41'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E26").Value = '  -1.87%  '
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '589.42'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E27").Value = '  +10.04%  '
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("E29").Value = '  +2.64%  '
$ws.Range("E30").Value = '  -1.32%  '
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("E33").Value = '  -1.19%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '6.61'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E34").Value = '  +3.77%  '
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '5.30'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E35").Value = '  -2.30%  '
$ws.Range("E36").Value = '  -1.36%  '
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '20.06'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E37").Value = '  -1.03%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E38").Value = '  +0.04%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '1.92'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E39").Value = '  +0.80%  '
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '151.68'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("E41").Value = '  -0.01%  '
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '159.80'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '2.38'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E43").Value = '  +3.88%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '4.02'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E44").Value = '  -0.78%  '
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '23.39'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E45").Value = '  +4.17%  '
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '0.0601'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("E47").Value = '  +4.47%  '
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '0.636'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E48").Value = '  +0.47%  '
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '0.0254'
$r.NumberFormat = "General"
$r.Style = "Normal"
$ws.Range("E49").Value = '  +0.32%  '
$ws.Range("D51").Value = '0.0₆0237'
$ws.Range("E51").Value = '  -6.50%  '
